# edit.ps1 - applies the two content edits described by the diff:
#   1. Slide 1, "CuadroTexto 3": "Rafael Gómez" -> "Rafael Gutiérrez"
#   2. Slide 4, "Marcador de contenido 2": rewords a sentence, splitting a
#      single run into five runs (same formatting, different text spans).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1: rename participant "Rafael Gómez" -> "Rafael Gutiérrez"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$participants = $slide1.Shapes.Item(3)   # "CuadroTexto 3"
$tr1 = $participants.TextFrame.TextRange
$para3 = $tr1.Paragraphs(3)              # "Rafael Gómez" paragraph
$run = $para3.Runs(1)
$run.Text = "Rafael Gutiérrez"

# ---------------------------------------------------------------------
# 2) Slide 4: reword "... donde almacenaremos la distancia de cada uno
#    de los satélites entre ellos. (Esto se realiza ..." into the new
#    phrasing, spread across five runs of identical formatting.
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$content = $slide4.Shapes.Item(2)        # "Marcador de contenido 2"
$tr2 = $content.TextFrame.TextRange
$para1 = $tr2.Paragraphs(1)

$targetRun = $para1.Runs(3)
$start = $targetRun.Start
$oldLen = $targetRun.Text.Length

$newParts = @(
    " donde almacenaremos la ",
    "distancia entre ",
    "cada uno de ",
    "los satélites. ",
    "(Esto se realiza en Matlab mediante el comando "
)

$newText = [string]::Join("", $newParts)

# Replace the whole old run's text with the new combined text in one go
# so the run's formatting (sz=2000, dirty="0") stays intact.
$whole = $tr2.Characters($start, $oldLen)
$whole.Text = $newText

# Re-split the combined text back into separate runs at each phrase
# boundary by re-assigning each span's Text to itself. This is a
# content no-op (text is unchanged) but makes the engine break out a
# fresh run with identical formatting at that exact span, matching the
# five-run structure without leaving any stray formatting attributes.
$offset = 0
foreach ($part in $newParts) {
    $segLen = $part.Length
    $seg = $tr2.Characters($start + $offset, $segLen)
    $seg.Text = $seg.Text
    $offset += $segLen
}
